$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the J column (k values)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17: labels in column A, aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the summary values bold/size-12, vertically centred, and raise the row height
$summaryValues = $ws.Range("B14:B17")
$summaryValues.Font.Bold = $true
$summaryValues.Font.Size = 12
$summaryValues.VerticalAlignment = -4108

$ws.Rows(14).RowHeight = 15.6
$ws.Rows(15).RowHeight = 15.6
$ws.Rows(16).RowHeight = 15.6
$ws.Rows(17).RowHeight = 15.6

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection matching the saved view state
$ws.Range("A14:B17").Select()
